$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert two new rows at the top of the data (below the header row) for
# the new daily entries (22 and 23 May), pushing the existing rows down.
$ws.Rows("2:3").Insert()

# Copy the date-format style from the row below (now row 4, the old row 2)
# onto the two new rows so they don't create a brand-new style entry.
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# New row for 23 May 2024 (serial 45435)
$ws.Range("A2").Value = 45435
$ws.Range("B2").Value = 243
$ws.Range("C2").Value = 26
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 210

# New row for 22 May 2024 (serial 45434)
$ws.Range("A3").Value = 45434
$ws.Range("B3").Value = 242
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 210

# Re-point the totals row formulas at the new data range (row 2 through
# the last data row, now row 24).
$ws.Range("B25").Formula = "=SUM(B2:B24)"
$ws.Range("C25:E25").Formula = "=SUM(C2:C24)"

# Notes pasted into column I alongside a few of the existing rows.
$ws.Range("I9").Value = "    timeCategories: {"
$ws.Range("I10").Value = "      'Regular arrivals': 210,"
$ws.Range("I8").Value = "count: 243,"
$ws.Range("I11").Value = "      'Night hour arrivals': 26,"
$ws.Range("I12").Value = "      'Shoulder hour flights': 7"

# Match the final selection state recorded in the workbook.
$ws.Range("B25:E25").Select()
